$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 700.38464
$ws.Range("I2").Value = 658.4286
$ws.Range("J2").Value = 749.3333
$ws.Range("K2").Value = 658.4286
$ws.Range("L2").Value = 749.3333
$ws.Range("M2").Value = -545.4286
$ws.Range("N2").Value = -975.3333

$ws.Range("H39").Value = 530.25
$ws.Range("I39").Value = 42
$ws.Range("J39").Value = 1018.5
$ws.Range("K39").Value = 126
$ws.Range("L39").Value = 3055.5
$ws.Range("M39").Value = 170
$ws.Range("N39").Value = -3647.5

$ws.Range("H40").Value = 2818.889
$ws.Range("J40").Value = 3396
$ws.Range("L40").Value = 3396
$ws.Range("N40").Value = -3746

$ws.Range("H111").Value = 4278.0435
$ws.Range("I111").Value = 4209.476
$ws.Range("K111").Value = 12628.428
$ws.Range("M111").Value = -9561.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2732.5
$ws.Range("I45").Value = 1706
$ws.Range("J45").Value = 3074.6667
$ws.Range("K45").Value = 1706
$ws.Range("L45").Value = 3074.6667
$ws.Range("M45").Value = -1329
$ws.Range("N45").Value = -3828.6667

$ws.Range("H63").Value = 7598.6924
$ws.Range("J63").Value = 9999.375
$ws.Range("L63").Value = 9999.375
$ws.Range("N63").Value = -11371.375

$ws.Range("H66").Value = 7598.6924
$ws.Range("J66").Value = 9999.375
$ws.Range("L66").Value = 49996.875
$ws.Range("N66").Value = -56860.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11549.72
$ws.Range("I86").Value = 3079.7144
$ws.Range("J86").Value = 22329.727
$ws.Range("K86").Value = 3079.7144
$ws.Range("L86").Value = 22329.727
$ws.Range("M86").Value = -1956.7144
$ws.Range("N86").Value = -24575.727

$ws.Range("H89").Value = 11549.72
$ws.Range("I89").Value = 3079.7144
$ws.Range("J89").Value = 22329.727
$ws.Range("K89").Value = 15398.572
$ws.Range("L89").Value = 111648.635
$ws.Range("M89").Value = -9782.572
$ws.Range("N89").Value = -122880.635

$ws.Range("H107").Value = 3251.077
$ws.Range("I107").Value = 2132.6155
$ws.Range("K107").Value = 2132.6155
$ws.Range("M107").Value = -212.6154999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7614.357
$ws.Range("I22").Value = 20248.2
$ws.Range("J22").Value = 595.55554
$ws.Range("K22").Value = 20248.2
$ws.Range("L22").Value = 595.55554
$ws.Range("M22").Value = -19898.2
$ws.Range("N22").Value = -1295.55554

$ws.Range("H53").Value = 165002
$ws.Range("J53").Value = 193752.75
$ws.Range("L53").Value = 193752.75
$ws.Range("N53").Value = -194966.75

$ws.Range("H58").Value = 2279.3447
$ws.Range("I58").Value = 1787.8182
$ws.Range("J58").Value = 3824.1428
$ws.Range("K58").Value = 1787.8182
$ws.Range("L58").Value = 3824.1428
$ws.Range("M58").Value = -1584.8182
$ws.Range("N58").Value = -4230.1428

$ws.Range("H74").Value = 44611.855
$ws.Range("I74").Value = 42285
$ws.Range("J74").Value = 44999.668
$ws.Range("K74").Value = 42285
$ws.Range("L74").Value = 44999.668
$ws.Range("M74").Value = -41411
$ws.Range("N74").Value = -46747.668

$ws.Range("H77").Value = 44611.855
$ws.Range("I77").Value = 42285
$ws.Range("J77").Value = 44999.668
$ws.Range("K77").Value = 126855
$ws.Range("L77").Value = 134999.004
$ws.Range("M77").Value = -122487
$ws.Range("N77").Value = -143735.004

$ws.Range("H136").Value = 2279.3447
$ws.Range("I136").Value = 1787.8182
$ws.Range("J136").Value = 3824.1428
$ws.Range("K136").Value = 5363.4546
$ws.Range("L136").Value = 11472.4284
$ws.Range("M136").Value = -2813.4546
$ws.Range("N136").Value = -16572.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 501.13333
$ws.Range("J12").Value = 525.9167
$ws.Range("L12").Value = 1577.7501
$ws.Range("N12").Value = -1923.7501

$ws.Range("H107").Value = 1095.625
$ws.Range("I107").Value = 224
$ws.Range("J107").Value = 1967.25
$ws.Range("K107").Value = 672
$ws.Range("L107").Value = 5901.75
$ws.Range("M107").Value = 1248
$ws.Range("N107").Value = -9741.75

$ws.Range("H113").Value = 2586.6191
$ws.Range("I113").Value = 1644.25
$ws.Range("K113").Value = 4932.75
$ws.Range("M113").Value = -2762.75

$ws.Range("H120").Value = 9716.166999999999
$ws.Range("I120").Value = 4446.6665
$ws.Range("J120").Value = 14985.667
$ws.Range("K120").Value = 13339.9995
$ws.Range("L120").Value = 44957.001
$ws.Range("M120").Value = -8501.999500000002
$ws.Range("N120").Value = -54633.001

$ws.Range("H122").Value = 1738.7
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws.Range("H129").Value = 4010.24
$ws.Range("I129").Value = 3893.7
$ws.Range("K129").Value = 11681.1
$ws.Range("M129").Value = -6681.099999999999

$ws.Range("H137").Value = 3701.125
$ws.Range("I137").Value = 2262.6667
$ws.Range("K137").Value = 6788.000100000001
$ws.Range("M137").Value = -1688.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3247.25
$ws.Range("I80").Value = 2849.7144
$ws.Range("J80").Value = 3803.8
$ws.Range("K80").Value = 2849.7144
$ws.Range("L80").Value = 3803.8
$ws.Range("M80").Value = -1851.7144
$ws.Range("N80").Value = -5799.8

$ws.Range("H83").Value = 3247.25
$ws.Range("I83").Value = 2849.7144
$ws.Range("J83").Value = 3803.8
$ws.Range("K83").Value = 14248.572
$ws.Range("L83").Value = 19019
$ws.Range("M83").Value = -9256.572
$ws.Range("N83").Value = -29003

$ws.Range("H107").Value = 399.38095
$ws.Range("I107").Value = 529.1818
$ws.Range("J107").Value = 256.6
$ws.Range("K107").Value = 529.1818
$ws.Range("L107").Value = 256.6
$ws.Range("M107").Value = 1390.8182
$ws.Range("N107").Value = -4096.6

$ws.Range("H113").Value = 2516.7188
$ws.Range("I113").Value = 1930.8823
$ws.Range("J113").Value = 3180.6667
$ws.Range("K113").Value = 1930.8823
$ws.Range("L113").Value = 3180.6667
$ws.Range("M113").Value = 239.1177
$ws.Range("N113").Value = -7520.6667

$ws.Range("H132").Value = 2092.8333
$ws.Range("I132").Value = 2133.9524
$ws.Range("J132").Value = 1805
$ws.Range("K132").Value = 6401.8572
$ws.Range("L132").Value = 5415
$ws.Range("M132").Value = -3871.8572
$ws.Range("N132").Value = -10475

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1278.2456
$ws.Range("J46").Value = 3327.2727
$ws.Range("L46").Value = 3327.2727
$ws.Range("N46").Value = -3703.2727

$ws.Range("H61").Value = 4250
$ws.Range("J61").Value = 6080
$ws.Range("L61").Value = 6080
$ws.Range("N61").Value = -6484

$ws.Range("H100").Value = 2814.2727
$ws.Range("I100").Value = 2183.8333
$ws.Range("K100").Value = 2183.8333
$ws.Range("M100").Value = -1642.8333

$ws.Range("H113").Value = 4250
$ws.Range("J113").Value = 6080
$ws.Range("L113").Value = 6080
$ws.Range("N113").Value = -10420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 34518570
$ws.Range("I122").Value = 38501270
$ws.Range("J122").Value = 1797
$ws.Range("K122").Value = 115503810
$ws.Range("L122").Value = 5391
$ws.Range("M122").Value = -115501360
$ws.Range("N122").Value = -10291

$ws.Range("H132").Value = 4309.9585
$ws.Range("I132").Value = 4553.525
$ws.Range("J132").Value = 3092.125
$ws.Range("K132").Value = 13660.575
$ws.Range("L132").Value = 9276.375
$ws.Range("M132").Value = -11130.575
$ws.Range("N132").Value = -14336.375
